# Update countries & provincias Spain
# Refreshed COVID-19 country dataset: re-ranked several countries by total
# cases (column B) and updated the case/recovery/death figures that moved
# with them, plus the "datos actualizados" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 16 de Junio de 2020 a las 22:17'

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = 'Estados Unidos'
$ws.Cells.Item(4, 2).Value = 2199303
$ws.Cells.Item(4, 3).Value = 16353
$ws.Cells.Item(4, 4).Value = 895168
$ws.Cells.Item(4, 5).Value = 1185241
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 611
$ws.Cells.Item(4, 8).Value = 118894

# Row 7: India
$ws.Cells.Item(7, 1).Value = 'India'
$ws.Cells.Item(7, 2).Value = 354148
$ws.Cells.Item(7, 3).Value = 11122
$ws.Cells.Item(7, 4).Value = 187552
$ws.Cells.Item(7, 5).Value = 154675
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 2006
$ws.Cells.Item(7, 8).Value = 11921

# Row 9: España
$ws.Cells.Item(9, 1).Value = 'España'
$ws.Cells.Item(9, 2).Value = 291408
$ws.Cells.Item(9, 3).Value = 219
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 27136

# Row 13: Alemania
$ws.Cells.Item(13, 1).Value = 'Alemania'
$ws.Cells.Item(13, 2).Value = 188373
$ws.Cells.Item(13, 3).Value = 329
$ws.Cells.Item(13, 4).Value = 173100
$ws.Cells.Item(13, 5).Value = 6364
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 24
$ws.Cells.Item(13, 8).Value = 8909

# Row 16: Francia
$ws.Cells.Item(16, 1).Value = 'Francia'
$ws.Cells.Item(16, 2).Value = 157716
$ws.Cells.Item(16, 3).Value = 344
$ws.Cells.Item(16, 4).Value = 73335
$ws.Cells.Item(16, 5).Value = 54834
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 111
$ws.Cells.Item(16, 8).Value = 29547

# Row 20: Canada
$ws.Cells.Item(20, 1).Value = 'Canada'
$ws.Cells.Item(20, 2).Value = 99426
$ws.Cells.Item(20, 3).Value = 279
$ws.Cells.Item(20, 4).Value = 61400
$ws.Cells.Item(20, 5).Value = 29813
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 38
$ws.Cells.Item(20, 8).Value = 8213

# Row 30: Egipto
$ws.Cells.Item(30, 1).Value = 'Egipto'
$ws.Cells.Item(30, 2).Value = 47856
$ws.Cells.Item(30, 3).Value = 1567
$ws.Cells.Item(30, 4).Value = 12730
$ws.Cells.Item(30, 5).Value = 33360
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 94
$ws.Cells.Item(30, 8).Value = 1766

# Row 31: Ecuador
$ws.Cells.Item(31, 1).Value = 'Ecuador'
$ws.Cells.Item(31, 2).Value = 47322
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 23349
$ws.Cells.Item(31, 5).Value = 20044
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 3929

# Row 49: Israel
$ws.Cells.Item(49, 1).Value = 'Israel'
$ws.Cells.Item(49, 2).Value = 19495
$ws.Cells.Item(49, 3).Value = 258
$ws.Cells.Item(49, 4).Value = 15449
$ws.Cells.Item(49, 5).Value = 3744
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 302

# Row 60: Ghana
$ws.Cells.Item(60, 1).Value = 'Ghana'
$ws.Cells.Item(60, 2).Value = 12193
$ws.Cells.Item(60, 3).Value = 229
$ws.Cells.Item(60, 4).Value = 4326
$ws.Cells.Item(60, 5).Value = 7809
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 4
$ws.Cells.Item(60, 8).Value = 58

# Row 61: Corea del Sur
$ws.Cells.Item(61, 1).Value = 'Corea del Sur'
$ws.Cells.Item(61, 2).Value = 12155
$ws.Cells.Item(61, 3).Value = 34
$ws.Cells.Item(61, 4).Value = 10760
$ws.Cells.Item(61, 5).Value = 1117
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 278

# Row 75: Costa de Marfil
$ws.Cells.Item(75, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(75, 2).Value = 5679
$ws.Cells.Item(75, 3).Value = 240
$ws.Cells.Item(75, 4).Value = 2637
$ws.Cells.Item(75, 5).Value = 2996
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 46

# Row 76: Uzbekistan
$ws.Cells.Item(76, 1).Value = 'Uzbekistan'
$ws.Cells.Item(76, 2).Value = 5493
$ws.Cells.Item(76, 3).Value = 230
$ws.Cells.Item(76, 4).Value = 4096
$ws.Cells.Item(76, 5).Value = 1378
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 19

# Row 80: Guinea
$ws.Cells.Item(80, 1).Value = 'Guinea'
$ws.Cells.Item(80, 2).Value = 4639
$ws.Cells.Item(80, 3).Value = 67
$ws.Cells.Item(80, 4).Value = 3327
$ws.Cells.Item(80, 5).Value = 1286
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 26

# Row 96: Kirguistan
$ws.Cells.Item(96, 1).Value = 'Kirguistan'
$ws.Cells.Item(96, 2).Value = 2472
$ws.Cells.Item(96, 3).Value = 100
$ws.Cells.Item(96, 4).Value = 1847
$ws.Cells.Item(96, 5).Value = 596
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 2
$ws.Cells.Item(96, 8).Value = 29

# Row 102: Mauritania
$ws.Cells.Item(102, 1).Value = 'Mauritania'
$ws.Cells.Item(102, 2).Value = 2057
$ws.Cells.Item(102, 3).Value = 170
$ws.Cells.Item(102, 4).Value = 373
$ws.Cells.Item(102, 5).Value = 1591
$ws.Cells.Item(102, 6).Value = 0
$ws.Cells.Item(102, 7).Value = 2
$ws.Cells.Item(102, 8).Value = 93

# Row 103: Estonia
$ws.Cells.Item(103, 1).Value = 'Estonia'
$ws.Cells.Item(103, 2).Value = 1975
$ws.Cells.Item(103, 3).Value = 1
$ws.Cells.Item(103, 4).Value = 1728
$ws.Cells.Item(103, 5).Value = 178
$ws.Cells.Item(103, 6).Value = 0
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 69

# Row 104: Sri Lanka
$ws.Cells.Item(104, 1).Value = 'Sri Lanka'
$ws.Cells.Item(104, 2).Value = 1915
$ws.Cells.Item(104, 3).Value = 10
$ws.Cells.Item(104, 4).Value = 1371
$ws.Cells.Item(104, 5).Value = 533
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 11

# Row 106: Nicaragua
$ws.Cells.Item(106, 1).Value = 'Nicaragua'
$ws.Cells.Item(106, 2).Value = 1823
$ws.Cells.Item(106, 3).Value = 359
$ws.Cells.Item(106, 4).Value = 1238
$ws.Cells.Item(106, 5).Value = 521
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 9
$ws.Cells.Item(106, 8).Value = 64

# Row 107: Islandia
$ws.Cells.Item(107, 1).Value = 'Islandia'
$ws.Cells.Item(107, 2).Value = 1812
$ws.Cells.Item(107, 3).Value = 2
$ws.Cells.Item(107, 4).Value = 1796
$ws.Cells.Item(107, 5).Value = 6
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 10

# Row 108: Costa Rica
$ws.Cells.Item(108, 1).Value = 'Costa Rica'
$ws.Cells.Item(108, 2).Value = 1796
$ws.Cells.Item(108, 3).Value = 52
$ws.Cells.Item(108, 4).Value = 794
$ws.Cells.Item(108, 5).Value = 990
$ws.Cells.Item(108, 6).Value = 0
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 12

# Row 109: Sudan del Sur
$ws.Cells.Item(109, 1).Value = 'Sudan del Sur'
$ws.Cells.Item(109, 2).Value = 1776
$ws.Cells.Item(109, 3).Value = 83
$ws.Cells.Item(109, 4).Value = 58
$ws.Cells.Item(109, 5).Value = 1688
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 3
$ws.Cells.Item(109, 8).Value = 30

# Row 110: Lituania
$ws.Cells.Item(110, 1).Value = 'Lituania'
$ws.Cells.Item(110, 2).Value = 1776
$ws.Cells.Item(110, 3).Value = 3
$ws.Cells.Item(110, 4).Value = 1441
$ws.Cells.Item(110, 5).Value = 259
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 76

# Row 111: Albania
$ws.Cells.Item(111, 1).Value = 'Albania'
$ws.Cells.Item(111, 2).Value = 1672
$ws.Cells.Item(111, 3).Value = 82
$ws.Cells.Item(111, 4).Value = 1064
$ws.Cells.Item(111, 5).Value = 571
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 1
$ws.Cells.Item(111, 8).Value = 37

# Row 112: Eslovaquia
$ws.Cells.Item(112, 1).Value = 'Eslovaquia'
$ws.Cells.Item(112, 2).Value = 1552
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(112, 4).Value = 1426
$ws.Cells.Item(112, 5).Value = 98
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 28

# Row 113: Nueva Zelanda
$ws.Cells.Item(113, 1).Value = 'Nueva Zelanda'
$ws.Cells.Item(113, 2).Value = 1506
$ws.Cells.Item(113, 3).Value = 2
$ws.Cells.Item(113, 4).Value = 1482
$ws.Cells.Item(113, 5).Value = 2
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 22

# Row 114: Eslovenia
$ws.Cells.Item(114, 1).Value = 'Eslovenia'
$ws.Cells.Item(114, 2).Value = 1499
$ws.Cells.Item(114, 3).Value = 3
$ws.Cells.Item(114, 4).Value = 1359
$ws.Cells.Item(114, 5).Value = 31
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 109

# Row 115: Guinea-Bisau
$ws.Cells.Item(115, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(115, 2).Value = 1492
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 153
$ws.Cells.Item(115, 5).Value = 1324
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 15

# Row 116: Libano
$ws.Cells.Item(116, 1).Value = 'Libano'
$ws.Cells.Item(116, 2).Value = 1473
$ws.Cells.Item(116, 3).Value = 9
$ws.Cells.Item(116, 4).Value = 889
$ws.Cells.Item(116, 5).Value = 552
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 32

# Row 117: Guayana Francesa
$ws.Cells.Item(117, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(117, 2).Value = 1421
$ws.Cells.Item(117, 3).Value = 95
$ws.Cells.Item(117, 4).Value = 619
$ws.Cells.Item(117, 5).Value = 797
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 2
$ws.Cells.Item(117, 8).Value = 5

# Row 118: Zambia
$ws.Cells.Item(118, 1).Value = 'Zambia'
$ws.Cells.Item(118, 2).Value = 1405
$ws.Cells.Item(118, 3).Value = 23
$ws.Cells.Item(118, 4).Value = 1142
$ws.Cells.Item(118, 5).Value = 252
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 11

# Row 121: Paraguay
$ws.Cells.Item(121, 1).Value = 'Paraguay'
$ws.Cells.Item(121, 2).Value = 1303
$ws.Cells.Item(121, 3).Value = 7
$ws.Cells.Item(121, 4).Value = 699
$ws.Cells.Item(121, 5).Value = 592
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 12

# Row 129: Burkina Faso
$ws.Cells.Item(129, 1).Value = 'Burkina Faso'
$ws.Cells.Item(129, 2).Value = 895
$ws.Cells.Item(129, 3).Value = 1
$ws.Cells.Item(129, 4).Value = 807
$ws.Cells.Item(129, 5).Value = 35
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 53

# Row 130: Yemen
$ws.Cells.Item(130, 1).Value = 'Yemen'
$ws.Cells.Item(130, 2).Value = 885
$ws.Cells.Item(130, 3).Value = 41
$ws.Cells.Item(130, 4).Value = 91
$ws.Cells.Item(130, 5).Value = 580
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 6
$ws.Cells.Item(130, 8).Value = 214

# Row 131: Congo
$ws.Cells.Item(131, 1).Value = 'Congo'
$ws.Cells.Item(131, 2).Value = 883
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 391
$ws.Cells.Item(131, 5).Value = 465
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 27

# Row 132: Georgia
$ws.Cells.Item(132, 1).Value = 'Georgia'
$ws.Cells.Item(132, 2).Value = 879
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 4).Value = 724
$ws.Cells.Item(132, 5).Value = 141
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 14

# Row 133: Principado de Andorra
$ws.Cells.Item(133, 1).Value = 'Principado de Andorra'
$ws.Cells.Item(133, 2).Value = 854
$ws.Cells.Item(133, 3).Value = 1
$ws.Cells.Item(133, 4).Value = 789
$ws.Cells.Item(133, 5).Value = 13
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 1
$ws.Cells.Item(133, 8).Value = 52

# Row 134: Republica del Chad
$ws.Cells.Item(134, 1).Value = 'Republica del Chad'
$ws.Cells.Item(134, 2).Value = 853
$ws.Cells.Item(134, 3).Value = 3
$ws.Cells.Item(134, 4).Value = 720
$ws.Cells.Item(134, 5).Value = 59
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 1
$ws.Cells.Item(134, 8).Value = 74

# Row 135: Uruguay
$ws.Cells.Item(135, 1).Value = 'Uruguay'
$ws.Cells.Item(135, 2).Value = 848
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = 792
$ws.Cells.Item(135, 5).Value = 33
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 23

# Row 136: Cabo Verde
$ws.Cells.Item(136, 1).Value = 'Cabo Verde'
$ws.Cells.Item(136, 2).Value = 781
$ws.Cells.Item(136, 3).Value = 21
$ws.Cells.Item(136, 4).Value = 354
$ws.Cells.Item(136, 5).Value = 420
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 7

# Row 140: Santo Tome y Principe
$ws.Cells.Item(140, 1).Value = 'Santo Tome y Principe'
$ws.Cells.Item(140, 2).Value = 671
$ws.Cells.Item(140, 3).Value = 9
$ws.Cells.Item(140, 4).Value = 182
$ws.Cells.Item(140, 5).Value = 477
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 12

# Row 143: Ruanda
$ws.Cells.Item(143, 1).Value = 'Ruanda'
$ws.Cells.Item(143, 2).Value = 636
$ws.Cells.Item(143, 3).Value = 24
$ws.Cells.Item(143, 4).Value = 338
$ws.Cells.Item(143, 5).Value = 296
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 2

# Row 144: Jamaica
$ws.Cells.Item(144, 1).Value = 'Jamaica'
$ws.Cells.Item(144, 2).Value = 621
$ws.Cells.Item(144, 3).Value = 4
$ws.Cells.Item(144, 4).Value = 430
$ws.Cells.Item(144, 5).Value = 181
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 10

# Row 145: Malaui
$ws.Cells.Item(145, 1).Value = 'Malaui'
$ws.Cells.Item(145, 2).Value = 564
$ws.Cells.Item(145, 3).Value = 9
$ws.Cells.Item(145, 4).Value = 73
$ws.Cells.Item(145, 5).Value = 485
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 6

# Row 146: Togo
$ws.Cells.Item(146, 1).Value = 'Togo'
$ws.Cells.Item(146, 2).Value = 537
$ws.Cells.Item(146, 3).Value = 6
$ws.Cells.Item(146, 4).Value = 344
$ws.Cells.Item(146, 5).Value = 180
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 13

# Row 147: Benin
$ws.Cells.Item(147, 1).Value = 'Benin'
$ws.Cells.Item(147, 2).Value = 532
$ws.Cells.Item(147, 3).Value = 49
$ws.Cells.Item(147, 4).Value = 236
$ws.Cells.Item(147, 5).Value = 287
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 9

# Row 155: Zimbabue
$ws.Cells.Item(155, 1).Value = 'Zimbabue'
$ws.Cells.Item(155, 2).Value = 391
$ws.Cells.Item(155, 3).Value = 4
$ws.Cells.Item(155, 4).Value = 62
$ws.Cells.Item(155, 5).Value = 325
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 4

# Row 161: Surinam
$ws.Cells.Item(161, 1).Value = 'Surinam'
$ws.Cells.Item(161, 2).Value = 236
$ws.Cells.Item(161, 3).Value = 7
$ws.Cells.Item(161, 4).Value = 48
$ws.Cells.Item(161, 5).Value = 182
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 1
$ws.Cells.Item(161, 8).Value = 6

# Row 171: Angola
$ws.Cells.Item(171, 1).Value = 'Angola'
$ws.Cells.Item(171, 2).Value = 148
$ws.Cells.Item(171, 3).Value = 6
$ws.Cells.Item(171, 4).Value = 64
$ws.Cells.Item(171, 5).Value = 78
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 0
$ws.Cells.Item(171, 8).Value = 6

# Row 172: Bermudas
$ws.Cells.Item(172, 1).Value = 'Bermudas'
$ws.Cells.Item(172, 2).Value = 144
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 4).Value = 127
$ws.Cells.Item(172, 5).Value = 8
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 9

# Row 180: Monaco
$ws.Cells.Item(180, 1).Value = 'Monaco'
$ws.Cells.Item(180, 2).Value = 99
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 94
$ws.Cells.Item(180, 5).Value = 1
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 4

# Row 191: Gambia
$ws.Cells.Item(191, 1).Value = 'Gambia'
$ws.Cells.Item(191, 2).Value = 34
$ws.Cells.Item(191, 3).Value = 4
$ws.Cells.Item(191, 4).Value = 24
$ws.Cells.Item(191, 5).Value = 9
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 1

# Row 192: Guam
$ws.Cells.Item(192, 1).Value = 'Guam'
$ws.Cells.Item(192, 2).Value = 32
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 5).Value = 31
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 1

